# Add a new "2022-Q4" sheet (copied from "2022-Q3") right after "总计",
# fill it with the new quarter's fund data, and record the new quarter
# in the "总计" summary sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "2022-Q4" sheet by copying "2022-Q3" -------------
# Copying before "2022-Q3" places the new sheet immediately after "总计".
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# --- 2. Fill in the 2022-Q4 fund data -----------------------------------
# Row 2: 002423 - 华宝标普美国品质消费股票（LOF）美元
$q4.Range("D2").Value = 3.86
$q4.Range("E2").Value = 94.64
$q4.Range("F2").Value = 2.88
$q4.Range("G2").Value = 0.1112

# Row 3: 162415 - 华宝标普美国品质消费股票（LOF）人民币A
$q4.Range("E3").Value = 94.64
$q4.Range("F3").Value = 2.88
$q4.Range("G3").Value = 0.0824

# Row 4: 009975 - 华宝标普美国品质消费股票（LOF）人民币C
$q4.Range("D4").Value = 1.00
$q4.Range("E4").Value = 94.64
$q4.Range("F4").Value = 2.88
$q4.Range("G4").Value = 0.0288

# --- 3. Update the "总计" summary sheet ----------------------------------
# Insert a new row so the existing quarters shift down, then record
# 2022-Q4's totals in the freed-up row.
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.22

Write-Output "2022-Q4 sheet added"
